$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115-209 down to 116-210
$ws.Rows(115).Insert()

# Populate the newly inserted row 115 with the new data record
$ws.Range("A115").Value = 9
$ws.Range("B115").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C115").Value = "Metropolitana"
$ws.Range("D115").Value = 44574
$ws.Range("E115").Value = 13
$ws.Range("F115").Value = 300000001
$ws.Range("G115").Value = "Rabanito"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 7900
$ws.Range("K115").Value = 2500
$ws.Range("L115").Value = 3000
$ws.Range("M115").Value = 2753
$ws.Range("N115").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O115").Value = "Provincia de Chacabuco"
$ws.Range("P115").Value = 28
$ws.Range("Q115").Value = 100
$ws.Range("R115").Value = "Hortaliza"
